$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

# nChildren
Set-TextValue "B4" "5"

# Row 6
Set-TextValue "B6" "6"
$ws.Range("C6").Value = "Ema  "
$ws.Range("D6").Value = "Ardell  "
$ws.Range("E6").Value = "4,-7"
$ws.Range("F6").Value = "Carley(grandmother): 0533587167"
Set-TextValue "H6" "17.0"

# Row 7
Set-TextValue "B7" "5"
$ws.Range("C7").Value = "Patti  "
$ws.Range("D7").Value = "Lavenia  "
$ws.Range("E7").Value = "5,-6"
$ws.Range("F7").Value = "Jennell(mother): 0503029941"
$ws.Range("G7").Value = "7:02:00"
Set-TextValue "H7" "15.0"

# Row 8
Set-TextValue "B8" "9"
$ws.Range("C8").Value = "Letha  "
$ws.Range("D8").Value = "Stephenie  "
$ws.Range("E8").Value = "5,-5"
$ws.Range("F8").Value = "Sibyl(mother): 0567328221"
$ws.Range("G8").Value = "7:04:00"
Set-TextValue "H8" "13.0"

# Row 9
Set-TextValue "B9" "7"
$ws.Range("C9").Value = "Wyatt  "
$ws.Range("D9").Value = "Willette  "
$ws.Range("E9").Value = "6,-4"
$ws.Range("F9").Value = "Antionette(father): 0557331799"
$ws.Range("G9").Value = "7:06:00"
Set-TextValue "H9" "11.0"

# Row 10
Set-TextValue "B10" "8"
$ws.Range("C10").Value = "Marni  "
$ws.Range("D10").Value = "Shanika  "
$ws.Range("E10").Value = "7,-4"
$ws.Range("F10").Value = "Lady(mother): 0560804012"
$ws.Range("G10").Value = "7:08:00"
Set-TextValue "H10" "9.0"

# Row 11 becomes "school" row (shifted up from row 14), H11 removed
$ws.Range("A11").Value = "school"
Set-TextValue "B11" "3"
$ws.Range("C11").Value = "Ironiah"
$ws.Range("D11").Value = "mySchool"
$ws.Range("E11").Value = "0,0"
$ws.Range("F11").Value = "Shir(secretary): 0523345098"
$ws.Range("G11").Value = "7:17:00"
$ws.Range("H11").ClearContents()

# Row 12 becomes "cost" row (shifted up from row 15), C12:H12 removed
$ws.Range("A12").Value = "cost"
Set-TextValue "B12" "39.0"
$ws.Range("C12:H12").ClearContents()

# Row 13 becomes "time" row (shifted up from row 16), C13:H13 removed
$ws.Range("A13").Value = "time"
Set-TextValue "B13" "17.0"
$ws.Range("C13:H13").ClearContents()

# Rows 14, 15, 16 are now entirely removed
$ws.Range("A14:H16").ClearContents()
